# Bugfixed the naive forecaster component module
#
# This script rewrites the dt_full_yoy_PUBCON_AR_50_9 worksheet data to reflect
# the corrected ("bugfixed") naive forecaster component forecasts:
#   - Row 1 (date headers) shifts: the first forecast-origin column (B) now
#     begins one period later, and values for every header cell are updated.
#   - Column BA (the last forecast-horizon column) is removed completely,
#     and rows 23:24 (the two latest forecast-origin rows) are removed
#     completely -> sheet dimension shrinks from A1:BA24 to A1:AZ22.
#   - For each remaining data row (3:22) the leading forecast cells that
#     preceded the corrected model's starting horizon are cleared, and the
#     remaining forecast cells are overwritten with the corrected values.
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Delete rows 23:24 entirely (they are removed completely in the target)
$ws.Rows("23:24").Delete()

# Step 2: Delete column BA entirely (removed completely in the target)
$ws.Columns("BA:BA").Delete()

# Row 1: header dates, full overwrite B1:AZ1 (values shifted/changed)
$row1 = New-Object 'object[,]' 1,51
$row1[0,0] = 39583
$row1[0,1] = 39765
$row1[0,2] = 39948
$row1[0,3] = 40130
$row1[0,4] = 40310
$row1[0,5] = 40494
$row1[0,6] = 40676
$row1[0,7] = 40862
$row1[0,8] = 41044
$row1[0,9] = 41228
$row1[0,10] = 41409
$row1[0,11] = 41592
$row1[0,12] = 41774
$row1[0,13] = 41957
$row1[0,14] = 42137
$row1[0,15] = 42321
$row1[0,16] = 42503
$row1[0,17] = 42689
$row1[0,18] = 42867
$row1[0,19] = 43053
$row1[0,20] = 43145
$row1[0,21] = 43235
$row1[0,22] = 43326
$row1[0,23] = 43418
$row1[0,24] = 43510
$row1[0,25] = 43600
$row1[0,26] = 43691
$row1[0,27] = 43783
$row1[0,28] = 43875
$row1[0,29] = 43966
$row1[0,30] = 44068
$row1[0,31] = 44159
$row1[0,32] = 44251
$row1[0,33] = 44341
$row1[0,34] = 44432
$row1[0,35] = 44525
$row1[0,36] = 44617
$row1[0,37] = 44706
$row1[0,38] = 44798
$row1[0,39] = 44890
$row1[0,40] = 44981
$row1[0,41] = 45071
$row1[0,42] = 45163
$row1[0,43] = 45254
$row1[0,44] = 45345
$row1[0,45] = 45436
$row1[0,46] = 45534
$row1[0,47] = 45618
$row1[0,48] = 45713
$row1[0,49] = 45800
$row1[0,50] = 45891
$ws.Range("B1:AZ1").Value = $row1

# Row 3
$ws.Range("C3:D3").ClearContents()
$row3 = New-Object 'object[,]' 1,48
$row3[0,0] = 1.834695583582535
$row3[0,1] = 1.834695583582535
$row3[0,2] = 1.834695583582535
$row3[0,3] = 1.834695583582535
$row3[0,4] = 1.834695583582535
$row3[0,5] = 1.834695583582535
$row3[0,6] = 1.834695583582491
$row3[0,7] = 1.834695583582491
$row3[0,8] = 1.834695583582491
$row3[0,9] = 1.834695583582491
$row3[0,10] = 1.834695583582491
$row3[0,11] = 1.834695583582491
$row3[0,12] = 1.834695583582491
$row3[0,13] = 1.834695583582491
$row3[0,14] = 1.834695583582491
$row3[0,15] = 1.834695583582491
$row3[0,16] = 1.834695583582491
$row3[0,17] = 1.834695583582491
$row3[0,18] = 1.834695583582491
$row3[0,19] = 1.834695583582491
$row3[0,20] = 1.834695583582491
$row3[0,21] = 1.834695583582491
$row3[0,22] = 1.834695583582491
$row3[0,23] = 1.834695583582491
$row3[0,24] = 1.834695583582491
$row3[0,25] = 1.834695583582491
$row3[0,26] = 1.834695583582491
$row3[0,27] = 1.834695583582491
$row3[0,28] = 1.834695583582491
$row3[0,29] = 1.834695583582491
$row3[0,30] = 1.834695583582491
$row3[0,31] = 1.834695583582491
$row3[0,32] = 1.834695583582491
$row3[0,33] = 1.834695583582491
$row3[0,34] = 1.834695583582491
$row3[0,35] = 1.834695583582491
$row3[0,36] = 1.834695583582491
$row3[0,37] = 1.834695583582491
$row3[0,38] = 1.834695583582491
$row3[0,39] = 1.834695583582491
$row3[0,40] = 1.834695583582491
$row3[0,41] = 1.834695583582491
$row3[0,42] = 1.834695583582491
$row3[0,43] = 1.834695583582491
$row3[0,44] = 1.834695583582491
$row3[0,45] = 1.834695583582491
$row3[0,46] = 1.834695583582491
$row3[0,47] = 1.834695583582491
$ws.Range("E3:AZ3").Value = $row3

# Row 4
$ws.Range("C4:F4").ClearContents()
$row4 = New-Object 'object[,]' 1,46
$row4[0,0] = 1.767835936772144
$row4[0,1] = 1.767835936772144
$row4[0,2] = 1.767835936772144
$row4[0,3] = 1.767835936772144
$row4[0,4] = 1.767835936772166
$row4[0,5] = 1.767835936772166
$row4[0,6] = 1.767835936772166
$row4[0,7] = 1.767835936772166
$row4[0,8] = 1.767835936772166
$row4[0,9] = 1.767835936772166
$row4[0,10] = 1.767835936772166
$row4[0,11] = 1.767835936772166
$row4[0,12] = 1.767835936772166
$row4[0,13] = 1.767835936772166
$row4[0,14] = 1.767835936772166
$row4[0,15] = 1.767835936772166
$row4[0,16] = 1.767835936772166
$row4[0,17] = 1.767835936772166
$row4[0,18] = 1.767835936772166
$row4[0,19] = 1.767835936772166
$row4[0,20] = 1.767835936772166
$row4[0,21] = 1.767835936772166
$row4[0,22] = 1.767835936772166
$row4[0,23] = 1.767835936772166
$row4[0,24] = 1.767835936772166
$row4[0,25] = 1.767835936772166
$row4[0,26] = 1.767835936772166
$row4[0,27] = 1.767835936772166
$row4[0,28] = 1.767835936772166
$row4[0,29] = 1.767835936772166
$row4[0,30] = 1.767835936772166
$row4[0,31] = 1.767835936772166
$row4[0,32] = 1.767835936772166
$row4[0,33] = 1.767835936772166
$row4[0,34] = 1.767835936772166
$row4[0,35] = 1.767835936772166
$row4[0,36] = 1.767835936772166
$row4[0,37] = 1.767835936772166
$row4[0,38] = 1.767835936772166
$row4[0,39] = 1.767835936772166
$row4[0,40] = 1.767835936772166
$row4[0,41] = 1.767835936772166
$row4[0,42] = 1.767835936772166
$row4[0,43] = 1.767835936772166
$row4[0,44] = 1.767835936772166
$row4[0,45] = 1.767835936772166
$ws.Range("G4:AZ4").Value = $row4

# Row 5
$ws.Range("C5:H5").ClearContents()
$row5 = New-Object 'object[,]' 1,44
$row5[0,0] = 1.074400434091038
$row5[0,1] = 1.074400434091038
$row5[0,2] = 1.074400434091016
$row5[0,3] = 1.074400434091016
$row5[0,4] = 1.074400434091016
$row5[0,5] = 1.074400434091016
$row5[0,6] = 1.074400434091016
$row5[0,7] = 1.074400434091016
$row5[0,8] = 1.074400434091016
$row5[0,9] = 1.074400434091016
$row5[0,10] = 1.074400434091016
$row5[0,11] = 1.074400434091016
$row5[0,12] = 1.074400434091016
$row5[0,13] = 1.074400434091016
$row5[0,14] = 1.074400434091016
$row5[0,15] = 1.074400434091016
$row5[0,16] = 1.074400434091016
$row5[0,17] = 1.074400434091016
$row5[0,18] = 1.074400434091016
$row5[0,19] = 1.074400434091016
$row5[0,20] = 1.074400434091016
$row5[0,21] = 1.074400434091016
$row5[0,22] = 1.074400434091016
$row5[0,23] = 1.074400434091016
$row5[0,24] = 1.074400434091016
$row5[0,25] = 1.074400434091016
$row5[0,26] = 1.074400434091016
$row5[0,27] = 1.074400434091016
$row5[0,28] = 1.074400434091016
$row5[0,29] = 1.074400434091016
$row5[0,30] = 1.074400434091016
$row5[0,31] = 1.074400434091016
$row5[0,32] = 1.074400434091016
$row5[0,33] = 1.074400434091016
$row5[0,34] = 1.074400434091016
$row5[0,35] = 1.074400434091016
$row5[0,36] = 1.074400434091016
$row5[0,37] = 1.074400434091016
$row5[0,38] = 1.074400434091016
$row5[0,39] = 1.074400434091016
$row5[0,40] = 1.074400434091016
$row5[0,41] = 1.074400434091016
$row5[0,42] = 1.074400434091016
$row5[0,43] = 1.074400434091016
$ws.Range("I5:AZ5").Value = $row5

# Row 6
$ws.Range("E6:J6").ClearContents()
$row6 = New-Object 'object[,]' 1,42
$row6[0,0] = 0.9212998022035679
$row6[0,1] = 0.9212998022035679
$row6[0,2] = 0.9212998022035679
$row6[0,3] = 0.9212998022035679
$row6[0,4] = 0.9212998022035679
$row6[0,5] = 0.9212998022035679
$row6[0,6] = 0.9212998022035679
$row6[0,7] = 0.9212998022035679
$row6[0,8] = 0.9212998022035679
$row6[0,9] = 0.9212998022035679
$row6[0,10] = 0.9212998022035679
$row6[0,11] = 0.9212998022035679
$row6[0,12] = 0.9212998022035679
$row6[0,13] = 0.9212998022035679
$row6[0,14] = 0.9212998022035679
$row6[0,15] = 0.9212998022035679
$row6[0,16] = 0.9212998022035679
$row6[0,17] = 0.9212998022035679
$row6[0,18] = 0.9212998022035679
$row6[0,19] = 0.9212998022035679
$row6[0,20] = 0.9212998022035679
$row6[0,21] = 0.9212998022035679
$row6[0,22] = 0.9212998022035679
$row6[0,23] = 0.9212998022035679
$row6[0,24] = 0.9212998022035679
$row6[0,25] = 0.9212998022035679
$row6[0,26] = 0.9212998022035679
$row6[0,27] = 0.9212998022035679
$row6[0,28] = 0.9212998022035679
$row6[0,29] = 0.9212998022035679
$row6[0,30] = 0.9212998022035679
$row6[0,31] = 0.9212998022035679
$row6[0,32] = 0.9212998022035679
$row6[0,33] = 0.9212998022035679
$row6[0,34] = 0.9212998022035679
$row6[0,35] = 0.9212998022035679
$row6[0,36] = 0.9212998022035679
$row6[0,37] = 0.9212998022035679
$row6[0,38] = 0.9212998022035679
$row6[0,39] = 0.9212998022035679
$row6[0,40] = 0.9212998022035679
$row6[0,41] = 0.9212998022035679
$ws.Range("K6:AZ6").Value = $row6

# Row 7
$ws.Range("G7:J7").ClearContents()
$row7 = New-Object 'object[,]' 1,42
$row7[0,0] = 1.274704633957136
$row7[0,1] = 1.097054137926201
$row7[0,2] = 1.141837882844188
$row7[0,3] = 1.141837882844188
$row7[0,4] = 1.141837882844188
$row7[0,5] = 1.141837882844188
$row7[0,6] = 1.141837882844188
$row7[0,7] = 1.141837882844188
$row7[0,8] = 1.141837882844188
$row7[0,9] = 1.141837882844188
$row7[0,10] = 1.141837882844188
$row7[0,11] = 1.141837882844188
$row7[0,12] = 1.141837882844188
$row7[0,13] = 1.141837882844188
$row7[0,14] = 1.141837882844188
$row7[0,15] = 1.141837882844188
$row7[0,16] = 1.141837882844188
$row7[0,17] = 1.141837882844188
$row7[0,18] = 1.141837882844188
$row7[0,19] = 1.141837882844188
$row7[0,20] = 1.141837882844188
$row7[0,21] = 1.141837882844188
$row7[0,22] = 1.141837882844188
$row7[0,23] = 1.141837882844188
$row7[0,24] = 1.141837882844188
$row7[0,25] = 1.141837882844188
$row7[0,26] = 1.141837882844188
$row7[0,27] = 1.141837882844188
$row7[0,28] = 1.141837882844188
$row7[0,29] = 1.141837882844188
$row7[0,30] = 1.141837882844188
$row7[0,31] = 1.141837882844188
$row7[0,32] = 1.141837882844188
$row7[0,33] = 1.141837882844188
$row7[0,34] = 1.141837882844188
$row7[0,35] = 1.141837882844188
$row7[0,36] = 1.141837882844188
$row7[0,37] = 1.141837882844188
$row7[0,38] = 1.141837882844188
$row7[0,39] = 1.141837882844188
$row7[0,40] = 1.141837882844188
$row7[0,41] = 1.141837882844188
$ws.Range("K7:AZ7").Value = $row7

# Row 8
$ws.Range("I8:J8").ClearContents()
$row8 = New-Object 'object[,]' 1,42
$row8[0,0] = 1.369652951216827
$row8[0,1] = 1.37755776875883
$row8[0,2] = 1.404348988410131
$row8[0,3] = 1.385527545913412
$row8[0,4] = 1.335361538769475
$row8[0,5] = 1.335361538769475
$row8[0,6] = 1.335361538769475
$row8[0,7] = 1.335361538769475
$row8[0,8] = 1.335361538769475
$row8[0,9] = 1.335361538769475
$row8[0,10] = 1.335361538769475
$row8[0,11] = 1.335361538769475
$row8[0,12] = 1.335361538769475
$row8[0,13] = 1.335361538769475
$row8[0,14] = 1.335361538769475
$row8[0,15] = 1.335361538769475
$row8[0,16] = 1.335361538769475
$row8[0,17] = 1.335361538769475
$row8[0,18] = 1.335361538769475
$row8[0,19] = 1.335361538769475
$row8[0,20] = 1.335361538769475
$row8[0,21] = 1.335361538769475
$row8[0,22] = 1.335361538769475
$row8[0,23] = 1.335361538769475
$row8[0,24] = 1.335361538769475
$row8[0,25] = 1.335361538769475
$row8[0,26] = 1.335361538769475
$row8[0,27] = 1.335361538769475
$row8[0,28] = 1.335361538769475
$row8[0,29] = 1.335361538769475
$row8[0,30] = 1.335361538769475
$row8[0,31] = 1.335361538769475
$row8[0,32] = 1.335361538769475
$row8[0,33] = 1.335361538769475
$row8[0,34] = 1.335361538769475
$row8[0,35] = 1.335361538769475
$row8[0,36] = 1.335361538769475
$row8[0,37] = 1.335361538769475
$row8[0,38] = 1.335361538769475
$row8[0,39] = 1.335361538769475
$row8[0,40] = 1.335361538769475
$row8[0,41] = 1.335361538769475
$ws.Range("K8:AZ8").Value = $row8

# Row 9
$ws.Range("K9").ClearContents()
$row9 = New-Object 'object[,]' 1,41
$row9[0,0] = 1.340476821534331
$row9[0,1] = 1.353022330671849
$row9[0,2] = 1.329814931661888
$row9[0,3] = 1.269653854937691
$row9[0,4] = 1.296301936385214
$row9[0,5] = 1.202048372526998
$row9[0,6] = 1.202048372526998
$row9[0,7] = 1.202048372526998
$row9[0,8] = 1.202048372526998
$row9[0,9] = 1.202048372526998
$row9[0,10] = 1.202048372526998
$row9[0,11] = 1.202048372526998
$row9[0,12] = 1.202048372526998
$row9[0,13] = 1.202048372526998
$row9[0,14] = 1.202048372526998
$row9[0,15] = 1.202048372526998
$row9[0,16] = 1.202048372526998
$row9[0,17] = 1.202048372526998
$row9[0,18] = 1.202048372526998
$row9[0,19] = 1.202048372526998
$row9[0,20] = 1.202048372526998
$row9[0,21] = 1.202048372526998
$row9[0,22] = 1.202048372526998
$row9[0,23] = 1.202048372526998
$row9[0,24] = 1.202048372526998
$row9[0,25] = 1.202048372526998
$row9[0,26] = 1.202048372526998
$row9[0,27] = 1.202048372526998
$row9[0,28] = 1.202048372526998
$row9[0,29] = 1.202048372526998
$row9[0,30] = 1.202048372526998
$row9[0,31] = 1.202048372526998
$row9[0,32] = 1.202048372526998
$row9[0,33] = 1.202048372526998
$row9[0,34] = 1.202048372526998
$row9[0,35] = 1.202048372526998
$row9[0,36] = 1.202048372526998
$row9[0,37] = 1.202048372526998
$row9[0,38] = 1.202048372526998
$row9[0,39] = 1.202048372526998
$row9[0,40] = 1.202048372526998
$ws.Range("L9:AZ9").Value = $row9

# Row 10
$ws.Range("M10").ClearContents()
$row10 = New-Object 'object[,]' 1,39
$row10[0,0] = 1.346941224959064
$row10[0,1] = 1.332063260659644
$row10[0,2] = 1.355477993452414
$row10[0,3] = 1.253742200752095
$row10[0,4] = 2.441628883342295
$row10[0,5] = 2.677488680362305
$row10[0,6] = 2.677488680362305
$row10[0,7] = 2.677488680362305
$row10[0,8] = 2.677488680362305
$row10[0,9] = 2.677488680362305
$row10[0,10] = 2.677488680362305
$row10[0,11] = 2.677488680362305
$row10[0,12] = 2.677488680362305
$row10[0,13] = 2.677488680362305
$row10[0,14] = 2.677488680362305
$row10[0,15] = 2.677488680362305
$row10[0,16] = 2.677488680362305
$row10[0,17] = 2.677488680362305
$row10[0,18] = 2.677488680362305
$row10[0,19] = 2.677488680362305
$row10[0,20] = 2.677488680362305
$row10[0,21] = 2.677488680362305
$row10[0,22] = 2.677488680362305
$row10[0,23] = 2.677488680362305
$row10[0,24] = 2.677488680362305
$row10[0,25] = 2.677488680362305
$row10[0,26] = 2.677488680362305
$row10[0,27] = 2.677488680362305
$row10[0,28] = 2.677488680362305
$row10[0,29] = 2.677488680362305
$row10[0,30] = 2.677488680362305
$row10[0,31] = 2.677488680362305
$row10[0,32] = 2.677488680362305
$row10[0,33] = 2.677488680362305
$row10[0,34] = 2.677488680362305
$row10[0,35] = 2.677488680362305
$row10[0,36] = 2.677488680362305
$row10[0,37] = 2.677488680362305
$row10[0,38] = 2.677488680362305
$ws.Range("N10:AZ10").Value = $row10

# Row 11
$ws.Range("O11").ClearContents()
$row11 = New-Object 'object[,]' 1,37
$row11[0,0] = 1.337624433084184
$row11[0,1] = 1.314367320499477
$row11[0,2] = 1.386547975635688
$row11[0,3] = 1.805615391969595
$row11[0,4] = 2.565764046666463
$row11[0,5] = 2.466954516646402
$row11[0,6] = 2.466954516646402
$row11[0,7] = 2.466954516646402
$row11[0,8] = 2.466954516646402
$row11[0,9] = 2.466954516646402
$row11[0,10] = 2.466954516646402
$row11[0,11] = 2.466954516646402
$row11[0,12] = 2.466954516646402
$row11[0,13] = 2.466954516646402
$row11[0,14] = 2.466954516646402
$row11[0,15] = 2.466954516646402
$row11[0,16] = 2.466954516646402
$row11[0,17] = 2.466954516646402
$row11[0,18] = 2.466954516646402
$row11[0,19] = 2.466954516646402
$row11[0,20] = 2.466954516646402
$row11[0,21] = 2.466954516646402
$row11[0,22] = 2.466954516646402
$row11[0,23] = 2.466954516646402
$row11[0,24] = 2.466954516646402
$row11[0,25] = 2.466954516646402
$row11[0,26] = 2.466954516646402
$row11[0,27] = 2.466954516646402
$row11[0,28] = 2.466954516646402
$row11[0,29] = 2.466954516646402
$row11[0,30] = 2.466954516646402
$row11[0,31] = 2.466954516646402
$row11[0,32] = 2.466954516646402
$row11[0,33] = 2.466954516646402
$row11[0,34] = 2.466954516646402
$row11[0,35] = 2.466954516646402
$row11[0,36] = 2.466954516646402
$ws.Range("P11:AZ11").Value = $row11

# Row 12
$ws.Range("Q12").ClearContents()
$row12 = New-Object 'object[,]' 1,35
$row12[0,0] = 1.458707167655282
$row12[0,1] = 1.519875552374694
$row12[0,2] = 1.833587970352424
$row12[0,3] = 1.661541796722577
$row12[0,4] = 1.376993627314671
$row12[0,5] = 1.263447557103259
$row12[0,6] = 1.362852986880547
$row12[0,7] = 1.401189216021326
$row12[0,8] = 1.401189216021326
$row12[0,9] = 1.401189216021326
$row12[0,10] = 1.401189216021326
$row12[0,11] = 1.401189216021326
$row12[0,12] = 1.401189216021326
$row12[0,13] = 1.401189216021326
$row12[0,14] = 1.401189216021326
$row12[0,15] = 1.401189216021326
$row12[0,16] = 1.401189216021326
$row12[0,17] = 1.401189216021326
$row12[0,18] = 1.401189216021326
$row12[0,19] = 1.401189216021326
$row12[0,20] = 1.401189216021326
$row12[0,21] = 1.401189216021326
$row12[0,22] = 1.401189216021326
$row12[0,23] = 1.401189216021326
$row12[0,24] = 1.401189216021326
$row12[0,25] = 1.401189216021326
$row12[0,26] = 1.401189216021326
$row12[0,27] = 1.401189216021326
$row12[0,28] = 1.401189216021326
$row12[0,29] = 1.401189216021326
$row12[0,30] = 1.401189216021326
$row12[0,31] = 1.401189216021326
$row12[0,32] = 1.401189216021326
$row12[0,33] = 1.401189216021326
$row12[0,34] = 1.401189216021326
$ws.Range("R12:AZ12").Value = $row12

# Row 13
$ws.Range("R13:S13").ClearContents()
$row13 = New-Object 'object[,]' 1,33
$row13[0,0] = 1.628446342551038
$row13[0,1] = 1.602599207906596
$row13[0,2] = 1.553425185274571
$row13[0,3] = 1.485511920344451
$row13[0,4] = 1.635045928803081
$row13[0,5] = 1.815016201748643
$row13[0,6] = 1.868356483387124
$row13[0,7] = 2.117022522597423
$row13[0,8] = 2.244229492550187
$row13[0,9] = 2.217567799050979
$row13[0,10] = 2.217567799050979
$row13[0,11] = 2.217567799050979
$row13[0,12] = 2.217567799050979
$row13[0,13] = 2.217567799050979
$row13[0,14] = 2.217567799050979
$row13[0,15] = 2.217567799050979
$row13[0,16] = 2.217567799050979
$row13[0,17] = 2.217567799050979
$row13[0,18] = 2.217567799050979
$row13[0,19] = 2.217567799050979
$row13[0,20] = 2.217567799050979
$row13[0,21] = 2.217567799050979
$row13[0,22] = 2.217567799050979
$row13[0,23] = 2.217567799050979
$row13[0,24] = 2.217567799050979
$row13[0,25] = 2.217567799050979
$row13[0,26] = 2.217567799050979
$row13[0,27] = 2.217567799050979
$row13[0,28] = 2.217567799050979
$row13[0,29] = 2.217567799050979
$row13[0,30] = 2.217567799050979
$row13[0,31] = 2.217567799050979
$row13[0,32] = 2.217567799050979
$ws.Range("T13:AZ13").Value = $row13

# Row 14
$ws.Range("T14:V14").ClearContents()
$row14 = New-Object 'object[,]' 1,30
$row14[0,0] = 1.556242995633506
$row14[0,1] = 1.576696032422364
$row14[0,2] = 1.59264109209869
$row14[0,3] = 1.600603483732033
$row14[0,4] = 1.745834498329324
$row14[0,5] = 1.944926006147352
$row14[0,6] = 1.810449264563152
$row14[0,7] = 2.005372766276703
$row14[0,8] = 2.149400276001101
$row14[0,9] = 2.139672475020404
$row14[0,10] = 2.139672475020404
$row14[0,11] = 2.139672475020404
$row14[0,12] = 2.139672475020404
$row14[0,13] = 2.139672475020404
$row14[0,14] = 2.139672475020404
$row14[0,15] = 2.139672475020404
$row14[0,16] = 2.139672475020404
$row14[0,17] = 2.139672475020404
$row14[0,18] = 2.139672475020404
$row14[0,19] = 2.139672475020404
$row14[0,20] = 2.139672475020404
$row14[0,21] = 2.139672475020404
$row14[0,22] = 2.139672475020404
$row14[0,23] = 2.139672475020404
$row14[0,24] = 2.139672475020404
$row14[0,25] = 2.139672475020404
$row14[0,26] = 2.139672475020404
$row14[0,27] = 2.139672475020404
$row14[0,28] = 2.139672475020404
$row14[0,29] = 2.139672475020404
$ws.Range("W14:AZ14").Value = $row14

# Row 15
$ws.Range("V15:Z15").ClearContents()
$row15 = New-Object 'object[,]' 1,26
$row15[0,0] = 1.634385169952313
$row15[0,1] = 1.659350806287785
$row15[0,2] = 1.646012149683163
$row15[0,3] = 1.679039075934385
$row15[0,4] = 1.76475225558832
$row15[0,5] = 1.674992401025355
$row15[0,6] = 2.128328071999674
$row15[0,7] = 2.273520076663971
$row15[0,8] = 2.453568910971748
$row15[0,9] = 2.100991693542231
$row15[0,10] = 2.100991693542231
$row15[0,11] = 2.100991693542231
$row15[0,12] = 2.100991693542231
$row15[0,13] = 2.100991693542231
$row15[0,14] = 2.100991693542231
$row15[0,15] = 2.100991693542231
$row15[0,16] = 2.100991693542231
$row15[0,17] = 2.100991693542231
$row15[0,18] = 2.100991693542231
$row15[0,19] = 2.100991693542231
$row15[0,20] = 2.100991693542231
$row15[0,21] = 2.100991693542231
$row15[0,22] = 2.100991693542231
$row15[0,23] = 2.100991693542231
$row15[0,24] = 2.100991693542231
$row15[0,25] = 2.100991693542231
$ws.Range("AA15:AZ15").Value = $row15

# Row 16
$ws.Range("Y16:AD16").ClearContents()
$row16 = New-Object 'object[,]' 1,22
$row16[0,0] = 1.683216193639248
$row16[0,1] = 1.676986128059976
$row16[0,2] = 1.787113358069314
$row16[0,3] = 1.827015901454287
$row16[0,4] = 2.131436976903012
$row16[0,5] = 0.5797569954799853
$row16[0,6] = 1.11435041103376
$row16[0,7] = 0.7380952286421882
$row16[0,8] = 0.812682184439506
$row16[0,9] = 0.8967077601845341
$row16[0,10] = 0.8967077601845341
$row16[0,11] = 0.8967077601845341
$row16[0,12] = 0.8967077601845341
$row16[0,13] = 0.8967077601845341
$row16[0,14] = 0.8967077601845341
$row16[0,15] = 0.8967077601845341
$row16[0,16] = 0.8967077601845341
$row16[0,17] = 0.8967077601845341
$row16[0,18] = 0.8967077601845341
$row16[0,19] = 0.8967077601845341
$row16[0,20] = 0.8967077601845341
$row16[0,21] = 0.8967077601845341
$ws.Range("AE16:AZ16").Value = $row16

# Row 17
$ws.Range("AC17:AG17").ClearContents()
$row17 = New-Object 'object[,]' 1,19
$row17[0,0] = 1.72853332973002
$row17[0,1] = 1.716088761125456
$row17[0,2] = 1.635325088332373
$row17[0,3] = 1.639897598917872
$row17[0,4] = 1.421947874507667
$row17[0,5] = 1.556352278772266
$row17[0,6] = 1.891749670939347
$row17[0,7] = 0.3338851812143995
$row17[0,8] = -0.4131792716363547
$row17[0,9] = 0.9940067218177528
$row17[0,10] = 0.782207885866093
$row17[0,11] = 0.782207885866093
$row17[0,12] = 0.782207885866093
$row17[0,13] = 0.782207885866093
$row17[0,14] = 0.782207885866093
$row17[0,15] = 0.782207885866093
$row17[0,16] = 0.782207885866093
$row17[0,17] = 0.782207885866093
$row17[0,18] = 0.782207885866093
$ws.Range("AH17:AZ17").Value = $row17

# Row 18
$ws.Range("AG18:AK18").ClearContents()
$row18 = New-Object 'object[,]' 1,15
$row18[0,0] = 1.586020868139215
$row18[0,1] = 1.63672950100644
$row18[0,2] = 1.696552096670256
$row18[0,3] = 1.536226507366623
$row18[0,4] = 1.173782198617435
$row18[0,5] = 1.820779918499094
$row18[0,6] = 1.655852464312013
$row18[0,7] = 2.228542839642689
$row18[0,8] = 1.805571054927801
$row18[0,9] = 1.634555928116921
$row18[0,10] = 1.508385007449875
$row18[0,11] = 1.508385007449875
$row18[0,12] = 1.508385007449875
$row18[0,13] = 1.508385007449875
$row18[0,14] = 1.508385007449875
$ws.Range("AL18:AZ18").Value = $row18

# Row 19
$ws.Range("AK19:AO19").ClearContents()
$row19 = New-Object 'object[,]' 1,11
$row19[0,0] = 1.487216163900373
$row19[0,1] = 1.708417111863691
$row19[0,2] = 1.700077181632453
$row19[0,3] = 1.712693234104812
$row19[0,4] = 1.632015075917925
$row19[0,5] = 1.554016159863814
$row19[0,6] = 1.310740881193517
$row19[0,7] = 0.9823016603409229
$row19[0,8] = 0.6395223689078522
$row19[0,9] = 0.6231570351797
$row19[0,10] = 0.6014263374495288
$ws.Range("AP19:AZ19").Value = $row19

# Row 20
$ws.Range("AO20:AS20").ClearContents()
$row20 = New-Object 'object[,]' 1,7
$row20[0,0] = 1.667247652000525
$row20[0,1] = 1.654472590788325
$row20[0,2] = 1.642837775344463
$row20[0,3] = 1.669808334748235
$row20[0,4] = 1.590970324046337
$row20[0,5] = 1.581524829939718
$row20[0,6] = 1.539180932382078
$ws.Range("AT20:AZ20").Value = $row20

# Row 21
$ws.Range("AS21:AW21").ClearContents()
$row21 = New-Object 'object[,]' 1,3
$row21[0,0] = 1.613721511063226
$row21[0,1] = 1.611425865651817
$row21[0,2] = 1.611459379842684
$ws.Range("AX21:AZ21").Value = $row21

# Row 22
$ws.Range("AW22:AZ22").ClearContents()

